$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "title"
$ws.Range("B3").Value = "Rock Quest"
$ws.Range("B3").Select()
